$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 20 onward down by 6 rows to make room for the
# new "Firm Size Definition" table (rows 21-25) above "Sector Distribution
# Details".
$ws.Range("A20:A25").EntireRow.Insert()

# New table header (row 21)
$ws.Range("B21").Value = "Number of employees"
$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B21:D21").Style = "title"

# Firm size definition rows (22-25)
$ws.Range("A22").Value = "Micro"
$ws.Range("B22").Value = "<5"

$ws.Range("A23").Value = "Small"
$ws.Range("B23").Value = "<50 Indstr. & Production<br/><25 Constr., Power engineering, Science, Education<br/><15 Transport, Trade, Services "

$ws.Range("A24").Value = "Medium"
$ws.Range("B24").Value = "<100 Indstr. & Production<br/><50 Constr., Power engineering, Science, Education<br/><30 Transport, Trade, Services "

$ws.Range("A25").Value = "Large"
$ws.Range("B25").Value = ">=100 Indstr. & Production<br/>>=50 Constr., Power engineering, Science, Education<br/>>=30 Transport, Trade, Services "

$ws.Range("A22:B25").Style = "Normal"

# Fix up the hyperlink, which the row insert does not auto-relocate: it
# used to sit on A36 and must now live on A42.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A42"), "http://www.unece.org/fileadmin/DAM/ceci/icp/Review/Studies/9.pdf")
$ws.Range("A42").Style = "HyperLink"
